$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Keep the Price column as text so numeric-looking values (e.g. "0.9998")
# are not auto-converted to numbers by Excel, matching the source data shape.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "31.104.98"
$ws.Range("E2").Value = "  +1.35%  "
$ws.Range("D3").Value = "1.956.02"
$ws.Range("E3").Value = "  +0.41%  "
$ws.Range("D4").Value = "0.9998"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "246.23"
$ws.Range("E6").Value = "  +0.05%  "
$ws.Range("D7").Value = "0.4907"
$ws.Range("E7").Value = "  +1.54%  "
$ws.Range("D8").Value = "0.2974"
$ws.Range("E8").Value = "  +1.14%  "
$ws.Range("D9").Value = "0.06858"
$ws.Range("E9").Value = "  +0.37%  "
$ws.Range("D10").Value = "19.17"
$ws.Range("E10").Value = "  -1.38%  "
$ws.Range("D11").Value = "108.08"
$ws.Range("E11").Value = "  -3.95%  "
$ws.Range("D12").Value = "0.07759"
$ws.Range("E12").Value = "  +1.26%  "
$ws.Range("D13").Value = "1.937.72"
$ws.Range("E13").Value = "  -0.74%  "
$ws.Range("D14").Value = "5.456"
$ws.Range("E14").Value = "  -1.28%  "
$ws.Range("D15").Value = "0.7083"
$ws.Range("E15").Value = "  +2.43%  "
$ws.Range("D16").Value = "284.49"
$ws.Range("E16").Value = "  -4.04%  "
$ws.Range("D17").Value = "31.115.23"
$ws.Range("E17").Value = "  +1.29%  "
$ws.Range("D18").Value = "0.000007774"
$ws.Range("E18").Value = "  +0.82%  "
$ws.Range("D19").Value = "13.25"
$ws.Range("E19").Value = "  -0.45%  "
$ws.Range("D20").Value = "1.000"
$ws.Range("E20").Value = "  +0.00%  "
$ws.Range("D21").Value = "2.184.50"
$ws.Range("E21").Value = "  -0.84%  "
$ws.Range("D22").Value = "5.499"
$ws.Range("E22").Value = "  -3.23%  "
$ws.Range("D23").Value = "1.000"
$ws.Range("E23").Value = "  -0.06%  "
$ws.Range("D24").Value = "6.519"
$ws.Range("E24").Value = "  -0.93%  "
$ws.Range("D25").Value = "9.818"
$ws.Range("E25").Value = "  -0.35%  "
$ws.Range("D26").Value = "169.18"
$ws.Range("E26").Value = "  +0.38%  "
$ws.Range("E27").Value = "  -1.34%  "
$ws.Range("D28").Value = "2.216"
$ws.Range("E28").Value = "  +1.04%  "
$ws.Range("D29").Value = "0.1057"
$ws.Range("E29").Value = "  -3.17%  "
$ws.Range("D30").Value = "1.425"
$ws.Range("E30").Value = "  -0.89%  "
$ws.Range("D31").Value = "1.584"
$ws.Range("E31").Value = "  -0.57%  "
$ws.Range("D32").Value = "4.577"
$ws.Range("E32").Value = "  -4.24%  "
$ws.Range("D33").Value = "4.456"
$ws.Range("E33").Value = "  -0.25%  "
$ws.Range("D34").Value = "0.04969"
$ws.Range("E34").Value = "  -2.60%  "
$ws.Range("D35").Value = "0.7575"
$ws.Range("E35").Value = "  -2.60%  "
$ws.Range("D36").Value = "1.182"
$ws.Range("E36").Value = "  +1.82%  "
$ws.Range("D37").Value = "2.733"
$ws.Range("E37").Value = "  -0.01%  "
$ws.Range("D38").Value = "0.02038"
$ws.Range("E38").Value = "  -2.13%  "
$ws.Range("D39").Value = "2.704"
$ws.Range("E39").Value = "  +0.16%  "
$ws.Range("D40").Value = "2.177"
$ws.Range("E40").Value = "  +5.93%  "
$ws.Range("D41").Value = "6.457"
$ws.Range("E41").Value = "  +8.93%  "
$ws.Range("B42").Value = "Aave"
$ws.Range("C42").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D42").Value = "74.13"
$ws.Range("E42").Value = "  +5.58%  "
$ws.Range("B43").Value = "TheSandbox"
$ws.Range("C43").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D43").Value = "0.4505"
$ws.Range("E43").Value = "  +0.55%  "
$ws.Range("D44").Value = "109.36"
$ws.Range("E44").Value = "  -1.59%  "
$ws.Range("D45").Value = "0.8825"
$ws.Range("E45").Value = "  +0.99%  "
$ws.Range("B46").Value = "Aptos"
$ws.Range("C46").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D46").Value = "8.176"
$ws.Range("E46").Value = "  +10.73%  "
$ws.Range("D47").Value = "1.001"
$ws.Range("E47").Value = "  -0.18%  "
$ws.Range("B48").Value = "Maker"
$ws.Range("C48").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D48").Value = "962.28"
$ws.Range("E48").Value = "  +5.82%  "
$ws.Range("B49").Value = "Algorand"
$ws.Range("C49").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D49").Value = "0.1267"
$ws.Range("E49").Value = "  +1.01%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "9.379"
$ws.Range("E50").Value = "  -0.86%  "
$ws.Range("D51").Value = "0.2587"
$ws.Range("E51").Value = "  +1.47%  "
